$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new weekly-report rows at the top of the data (row 2),
# pushing the existing rows (old 2-6) down to become rows 4-8.
$ws.Rows.Item(2).Resize(2).Insert()

# The insert copies formatting from the row above (the bold header row).
# Clear that out so the new rows match the plain data-row formatting,
# then restore the date number format on column D (same as other rows).
$ws.Range("A2:T3").ClearFormats()
$ws.Range("D2:D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Common values shared by the two new rows
foreach ($r in 2,3) {
    $ws.Range("A$r").Value = 11
    $ws.Range("B$r").Value = "Vega Monumental Concepción"
    $ws.Range("C$r").Value = "Bíobío"
    $ws.Range("E$r").Value = 8
    $ws.Range("F$r").Value = "Fruta"
    $ws.Range("G$r").Value = 100107
    $ws.Range("H$r").Value = "Otros"
    $ws.Range("I$r").Value = 100107011
    $ws.Range("J$r").Value = "Tuna"
    $ws.Range("K$r").Value = "Sin especificar"
    $ws.Range("Q$r").Value = "`$/caja 18 kilos"
    $ws.Range("R$r").Value = "Provincia de Melipilla"
    $ws.Range("T$r").Value = 18
}

# New row 2 (Primera, week of 44699)
$ws.Range("D2").Value = 44699
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 21000
$ws.Range("S2").Value = 1167

# New row 3 (Segunda, week of 44699)
$ws.Range("D3").Value = 44699
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("S3").Value = 1000
